$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.1243246666666667
$ws.Range("H2").Value = 0.372974
$ws.Range("I2").Value = 0.09963085929726231
$ws.Range("J2").Value = 0.09963085929726233
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.854571666666667
$ws.Range("N2").Value = 5.563715
$ws.Range("O2").Value = 0.01651371646154392
$ws.Range("P2").Value = 0.01651371646154392
$ws.Range("Q2").Value = 0.2305690042677778
$ws.Range("R2").Value = 2.07512103841
$ws.Range("S2").Value = 0.001645275761254967
$ws.Range("T2").Value = 0.001645275761254967
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.1243246666666667
$ws.Range("H3").Value = 0.372974
$ws.Range("I3").Value = 0.09963085929726231
$ws.Range("J3").Value = 0.09963085929726233
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 92.91372433333333
$ws.Range("N3").Value = 278.741173
$ws.Range("O3").Value = 0.8273343794712995
$ws.Range("P3").Value = 0.8273343794712996
$ws.Range("Q3").Value = 11.55146780650022
$ws.Range("R3").Value = 103.963210258502
$ws.Range("S3").Value = 0.08242803515289286
$ws.Range("T3").Value = 0.08242803515289289
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.1243246666666667
$ws.Range("H4").Value = 0.372974
$ws.Range("I4").Value = 0.09963085929726231
$ws.Range("J4").Value = 0.09963085929726233
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.162136
$ws.Range("N4").Value = 0.4864080000000001
$ws.Range("O4").Value = 0.001443712303133186
$ws.Range("P4").Value = 0.001443712303133187
$ws.Range("Q4").Value = 0.02015750415466667
$ws.Range("R4").Value = 0.181417537392
$ws.Range("S4").Value = 0.000143838297339189
$ws.Range("T4").Value = 0.0001438382973391891
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.1243246666666667
$ws.Range("H5").Value = 0.372974
$ws.Range("I5").Value = 0.09963085929726231
$ws.Range("J5").Value = 0.09963085929726233
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.37449166666667
$ws.Range("N5").Value = 52.123475
$ws.Range("O5").Value = 0.1547081917640233
$ws.Range("P5").Value = 0.1547081917640233
$ws.Range("Q5").Value = 2.160077884961111
$ws.Range("R5").Value = 19.44070096465
$ws.Range("S5").Value = 0.01541371008577528
$ws.Range("T5").Value = 0.01541371008577529
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.7328223333333334
$ws.Range("H6").Value = 2.198467
$ws.Range("I6").Value = 0.5872665557027417
$ws.Range("J6").Value = 0.5872665557027417
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.854571666666667
$ws.Range("N6").Value = 5.563715
$ws.Range("O6").Value = 0.01651371646154392
$ws.Range("P6").Value = 0.01651371646154392
$ws.Range("Q6").Value = 1.359071536100556
$ws.Range("R6").Value = 12.231643824905
$ws.Range("S6").Value = 0.009697953388222564
$ws.Range("T6").Value = 0.009697953388222566
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.7328223333333334
$ws.Range("H7").Value = 2.198467
$ws.Range("I7").Value = 0.5872665557027417
$ws.Range("J7").Value = 0.5872665557027417
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 92.91372433333333
$ws.Range("N7").Value = 278.741173
$ws.Range("O7").Value = 0.8273343794712995
$ws.Range("P7").Value = 0.8273343794712996
$ws.Range("Q7").Value = 68.08925226464345
$ws.Range("R7").Value = 612.803270381791
$ws.Range("S7").Value = 0.4858658114465751
$ws.Range("T7").Value = 0.4858658114465752
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.7328223333333334
$ws.Range("H8").Value = 2.198467
$ws.Range("I8").Value = 0.5872665557027417
$ws.Range("J8").Value = 0.5872665557027417
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.162136
$ws.Range("N8").Value = 0.4864080000000001
$ws.Range("O8").Value = 0.001443712303133186
$ws.Range("P8").Value = 0.001443712303133187
$ws.Range("Q8").Value = 0.1188168818373334
$ws.Range("R8").Value = 1.069351936536
$ws.Range("S8").Value = 0.0008478439516866988
$ws.Range("T8").Value = 0.000847843951686699
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.7328223333333334
$ws.Range("H9").Value = 2.198467
$ws.Range("I9").Value = 0.5872665557027417
$ws.Range("J9").Value = 0.5872665557027417
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 17.37449166666667
$ws.Range("N9").Value = 52.123475
$ws.Range("O9").Value = 0.1547081917640233
$ws.Range("P9").Value = 0.1547081917640233
$ws.Range("Q9").Value = 12.73241552364722
$ws.Range("R9").Value = 114.591739712825
$ws.Range("S9").Value = 0.09085494691625723
$ws.Range("T9").Value = 0.09085494691625723
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.2213483333333333
$ws.Range("H10").Value = 0.664045
$ws.Range("I10").Value = 0.1773833402919521
$ws.Range("J10").Value = 0.1773833402919521
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.854571666666667
$ws.Range("N10").Value = 5.563715
$ws.Range("O10").Value = 0.01651371646154392
$ws.Range("P10").Value = 0.01651371646154392
$ws.Range("Q10").Value = 0.4105063474638889
$ws.Range("R10").Value = 3.694557127175
$ws.Range("S10").Value = 0.002929258186582856
$ws.Range("T10").Value = 0.002929258186582857
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.2213483333333333
$ws.Range("H11").Value = 0.664045
$ws.Range("I11").Value = 0.1773833402919521
$ws.Range("J11").Value = 0.1773833402919521
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 92.91372433333333
$ws.Range("N11").Value = 278.741173
$ws.Range("O11").Value = 0.8273343794712995
$ws.Range("P11").Value = 0.8273343794712996
$ws.Range("Q11").Value = 20.56629802497611
$ws.Range("R11").Value = 185.096682224785
$ws.Range("S11").Value = 0.1467553357689886
$ws.Range("T11").Value = 0.1467553357689886
# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.2213483333333333
$ws.Range("H12").Value = 0.664045
$ws.Range("I12").Value = 0.1773833402919521
$ws.Range("J12").Value = 0.1773833402919521
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.162136
$ws.Range("N12").Value = 0.4864080000000001
$ws.Range("O12").Value = 0.001443712303133186
$ws.Range("P12").Value = 0.001443712303133187
$ws.Range("Q12").Value = 0.03588853337333334
$ws.Range("R12").Value = 0.3229968003600001
$ws.Range("S12").Value = 0.0002560905107503519
$ws.Range("T12").Value = 0.0002560905107503519
# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.2213483333333333
$ws.Range("H13").Value = 0.664045
$ws.Range("I13").Value = 0.1773833402919521
$ws.Range("J13").Value = 0.1773833402919521
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 17.37449166666667
$ws.Range("N13").Value = 52.123475
$ws.Range("O13").Value = 0.1547081917640233
$ws.Range("P13").Value = 0.1547081917640233
$ws.Range("Q13").Value = 3.845814772930556
$ws.Range("R13").Value = 34.612332956375
$ws.Range("S13").Value = 0.02744265582563033
$ws.Range("T13").Value = 0.02744265582563033
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.1693576666666667
$ws.Range("H14").Value = 0.508073
$ws.Range("I14").Value = 0.1357192447080439
$ws.Range("J14").Value = 0.1357192447080439
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.854571666666667
$ws.Range("N14").Value = 5.563715
$ws.Range("O14").Value = 0.01651371646154392
$ws.Range("P14").Value = 0.01651371646154392
$ws.Range("Q14").Value = 0.3140859301327777
$ws.Range("R14").Value = 2.826773371195
$ws.Range("S14").Value = 0.002241229125483532
$ws.Range("T14").Value = 0.002241229125483532
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.1693576666666667
$ws.Range("H15").Value = 0.508073
$ws.Range("I15").Value = 0.1357192447080439
$ws.Range("J15").Value = 0.1357192447080439
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 92.91372433333333
$ws.Range("N15").Value = 278.741173
$ws.Range("O15").Value = 0.8273343794712995
$ws.Range("P15").Value = 0.8273343794712996
$ws.Range("Q15").Value = 15.73565155440322
$ws.Range("R15").Value = 141.620863989629
$ws.Range("S15").Value = 0.1122851971028429
$ws.Range("T15").Value = 0.1122851971028429
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.1693576666666667
$ws.Range("H16").Value = 0.508073
$ws.Range("I16").Value = 0.1357192447080439
$ws.Range("J16").Value = 0.1357192447080439
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.162136
$ws.Range("N16").Value = 0.4864080000000001
$ws.Range("O16").Value = 0.001443712303133186
$ws.Range("P16").Value = 0.001443712303133187
$ws.Range("Q16").Value = 0.02745897464266667
$ws.Range("R16").Value = 0.247130771784
$ws.Range("S16").Value = 0.0001959395433569465
$ws.Range("T16").Value = 0.0001959395433569465
# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.1693576666666667
$ws.Range("H17").Value = 0.508073
$ws.Range("I17").Value = 0.1357192447080439
$ws.Range("J17").Value = 0.1357192447080439
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 17.37449166666667
$ws.Range("N17").Value = 52.123475
$ws.Range("O17").Value = 0.1547081917640233
$ws.Range("P17").Value = 0.1547081917640233
$ws.Range("Q17").Value = 2.942503368186111
$ws.Range("R17").Value = 26.482530313675
$ws.Range("S17").Value = 0.02099687893636045
$ws.Range("T17").Value = 0.02099687893636045
